# Insert two new data rows at 201-202 (pushes the former rows 201..310 down to 203..312)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("201:202").Insert()

# --- New row 201 ---
$ws.Range("A201").Value = 10
$ws.Range("B201").Value = "Vega Modelo de Temuco"
$ws.Range("C201").Value = "La Araucanía"
$ws.Range("D201").Value = 44992
$ws.Range("E201").Value = 9
$ws.Range("F201").Value = 100112043
$ws.Range("G201").Value = "Pepino dulce"
$ws.Range("H201").Value = "Cultivar XV región"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 115
$ws.Range("K201").Value = 17000
$ws.Range("L201").Value = 18000
$ws.Range("M201").Value = 17565
$ws.Range("N201").Value = "$/bandeja 18 kilos"
$ws.Range("O201").Value = "Región de Arica y Parinacota"
$ws.Range("P201").Value = 976
$ws.Range("Q201").Value = 18
$ws.Range("R201").Value = "Hortaliza"

# --- New row 202 ---
$ws.Range("A202").Value = 10
$ws.Range("B202").Value = "Vega Modelo de Temuco"
$ws.Range("C202").Value = "La Araucanía"
$ws.Range("D202").Value = 44992
$ws.Range("E202").Value = 9
$ws.Range("F202").Value = 100112043
$ws.Range("G202").Value = "Pepino dulce"
$ws.Range("H202").Value = "Cultivar XV región"
$ws.Range("I202").Value = "Segunda"
$ws.Range("J202").Value = 25
$ws.Range("K202").Value = 10000
$ws.Range("L202").Value = 10000
$ws.Range("M202").Value = 10000
$ws.Range("N202").Value = "$/bandeja 18 kilos"
$ws.Range("O202").Value = "Región de Arica y Parinacota"
$ws.Range("P202").Value = 556
$ws.Range("Q202").Value = 18
$ws.Range("R202").Value = "Hortaliza"
